$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1612
$ws.Range("I70").Value = 1998.5
$ws.Range("J70").Value = 1526.1111
$ws.Range("K70").Value = 5995.5
$ws.Range("L70").Value = 4578.3333
$ws.Range("M70").Value = -5725.5
$ws.Range("N70").Value = -5118.3333

$ws.Range("H73").Value = 1612
$ws.Range("I73").Value = 1998.5
$ws.Range("J73").Value = 1526.1111
$ws.Range("K73").Value = 5995.5
$ws.Range("L73").Value = 4578.3333
$ws.Range("M73").Value = -5059.5
$ws.Range("N73").Value = -6450.3333

$ws.Range("H106").Value = 85693.875
$ws.Range("I106").Value = 94221.57000000001
$ws.Range("J106").Value = 26000
$ws.Range("K106").Value = 94221.57000000001
$ws.Range("L106").Value = 26000
$ws.Range("M106").Value = -93590.57000000001
$ws.Range("N106").Value = -27262

$ws.Range("H118").Value = 497.8
$ws.Range("I118").Value = 520.3333
$ws.Range("K118").Value = 1560.9999
$ws.Range("M118").Value = 96.00009999999997

$ws.Range("H132").Value = 2053.4888
$ws.Range("I132").Value = 1945.2
$ws.Range("K132").Value = 5835.6
$ws.Range("M132").Value = -3305.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2204.35
$ws.Range("I45").Value = 2209.7
$ws.Range("J45").Value = 2199
$ws.Range("K45").Value = 2209.7
$ws.Range("L45").Value = 2199
$ws.Range("M45").Value = -1832.7
$ws.Range("N45").Value = -2953

$ws.Range("H88").Value = 1734.6072
$ws.Range("I88").Value = 1262.3334
$ws.Range("J88").Value = 2088.8125
$ws.Range("K88").Value = 1262.3334
$ws.Range("L88").Value = 2088.8125
$ws.Range("M88").Value = -856.3334
$ws.Range("N88").Value = -2900.8125

$ws.Range("H91").Value = 1734.6072
$ws.Range("I91").Value = 1262.3334
$ws.Range("J91").Value = 2088.8125
$ws.Range("K91").Value = 1262.3334
$ws.Range("L91").Value = 2088.8125
$ws.Range("M91").Value = 141.6666
$ws.Range("N91").Value = -4896.8125

$ws.Range("H97").Value = 324.9
$ws.Range("I97").Value = 299
$ws.Range("J97").Value = 402.6
$ws.Range("K97").Value = 299
$ws.Range("L97").Value = 402.6
$ws.Range("M97").Value = 197
$ws.Range("N97").Value = -1394.6

$ws.Range("H122").Value = 3081.8838
$ws.Range("I122").Value = 2858.0334
$ws.Range("J122").Value = 3598.4614
$ws.Range("K122").Value = 8574.100199999999
$ws.Range("L122").Value = 10795.3842
$ws.Range("M122").Value = -6124.100199999999
$ws.Range("N122").Value = -15695.3842

$ws.Range("H132").Value = 7388.8047
$ws.Range("I132").Value = 4949.387
$ws.Range("J132").Value = 14951
$ws.Range("K132").Value = 14848.161
$ws.Range("L132").Value = 44853
$ws.Range("M132").Value = -12318.161
$ws.Range("N132").Value = -49913

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2873.375
$ws.Range("I16").Value = 4999
$ws.Range("J16").Value = 2569.7144
$ws.Range("K16").Value = 4999
$ws.Range("L16").Value = 2569.7144
$ws.Range("M16").Value = -4712
$ws.Range("N16").Value = -3143.7144

$ws.Range("H22").Value = 1046.44
$ws.Range("I22").Value = 308.25
$ws.Range("J22").Value = 2358.7778
$ws.Range("K22").Value = 308.25
$ws.Range("L22").Value = 2358.7778
$ws.Range("M22").Value = 41.75
$ws.Range("N22").Value = -3058.7778

$ws.Range("H62").Value = 6076.857
$ws.Range("J62").Value = 5389
$ws.Range("L62").Value = 5389
$ws.Range("N62").Value = -6637

$ws.Range("H65").Value = 6076.857
$ws.Range("J65").Value = 5389
$ws.Range("L65").Value = 26945
$ws.Range("N65").Value = -33185

$ws.Range("H86").Value = 111131750
$ws.Range("I86").Value = 200015790
$ws.Range("J86").Value = 26697.5
$ws.Range("K86").Value = 200015790
$ws.Range("L86").Value = 26697.5
$ws.Range("M86").Value = -200014667
$ws.Range("N86").Value = -28943.5

$ws.Range("H89").Value = 111131750
$ws.Range("I89").Value = 200015790
$ws.Range("J89").Value = 26697.5
$ws.Range("K89").Value = 1000078950
$ws.Range("L89").Value = 133487.5
$ws.Range("M89").Value = -1000073334
$ws.Range("N89").Value = -144719.5

$ws.Range("H99").Value = 3207.0356
$ws.Range("I99").Value = 2788.4783
$ws.Range("J99").Value = 5132.4
$ws.Range("K99").Value = 2788.4783
$ws.Range("L99").Value = 5132.4
$ws.Range("M99").Value = -1290.4783
$ws.Range("N99").Value = -8128.4

$ws.Range("H107").Value = 1063.7931
$ws.Range("I107").Value = 651
$ws.Range("J107").Value = 1171.4783
$ws.Range("K107").Value = 651
$ws.Range("L107").Value = 1171.4783
$ws.Range("M107").Value = 1269
$ws.Range("N107").Value = -5011.4783

$ws.Range("H113").Value = 2873.375
$ws.Range("I113").Value = 4999
$ws.Range("J113").Value = 2569.7144
$ws.Range("K113").Value = 4999
$ws.Range("L113").Value = 2569.7144
$ws.Range("M113").Value = -2829
$ws.Range("N113").Value = -6909.7144

$ws.Range("H126").Value = 3207.0356
$ws.Range("I126").Value = 2788.4783
$ws.Range("J126").Value = 5132.4
$ws.Range("K126").Value = 8365.4349
$ws.Range("L126").Value = 15397.2
$ws.Range("M126").Value = -5895.4349
$ws.Range("N126").Value = -20337.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 2015.5625
$ws.Range("I59").Value = 3374.5
$ws.Range("J59").Value = 1821.4286
$ws.Range("K59").Value = 10123.5
$ws.Range("L59").Value = 5464.2858
$ws.Range("M59").Value = -9583.5
$ws.Range("N59").Value = -6544.2858

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H68").Value = 851.1667
$ws.Range("I68").Value = 600
$ws.Range("J68").Value = 901.4
$ws.Range("K68").Value = 1800
$ws.Range("L68").Value = 2704.2
$ws.Range("M68").Value = -989
$ws.Range("N68").Value = -4326.2

$ws.Range("H69").Value = 7500
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H71").Value = 851.1667
$ws.Range("I71").Value = 600
$ws.Range("J71").Value = 901.4
$ws.Range("K71").Value = 5400
$ws.Range("L71").Value = 8112.599999999999
$ws.Range("M71").Value = -1344
$ws.Range("N71").Value = -16224.6

$ws.Range("H72").Value = 7500
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3420
$ws.Range("I102").Value = 3228.125
$ws.Range("J102").Value = 4187.5
$ws.Range("K102").Value = 3228.125
$ws.Range("L102").Value = 4187.5
$ws.Range("M102").Value = -1606.125
$ws.Range("N102").Value = -7431.5

$ws.Range("H132").Value = 3332.5
$ws.Range("I132").Value = 3332.5
$ws.Range("K132").Value = 9997.5
$ws.Range("M132").Value = -7467.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1643.8
$ws.Range("I16").Value = 1679.75
$ws.Range("K16").Value = 1679.75
$ws.Range("M16").Value = -1509.75

$ws.Range("H98").Value = 355
$ws.Range("J98").Value = 355
$ws.Range("L98").Value = 355
$ws.Range("N98").Value = -6345

$ws.Range("H132").Value = 3473.4333
$ws.Range("I132").Value = 3235.5217
$ws.Range("J132").Value = 4255.143
$ws.Range("K132").Value = 9706.5651
$ws.Range("L132").Value = 12765.429
$ws.Range("M132").Value = -7176.5651
$ws.Range("N132").Value = -17825.429

$ws.Range("H138").Value = 99998.86
$ws.Range("J138").Value = 99998.86
$ws.Range("L138").Value = 99998.86
$ws.Range("N138").Value = -110278.86

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 9376.143
$ws.Range("I41").Value = 11436.5
$ws.Range("K41").Value = 11436.5
$ws.Range("M41").Value = -11046.5

$ws.Range("H132").Value = 1379.091
$ws.Range("I132").Value = 1501.625
$ws.Range("J132").Value = 1052.3334
$ws.Range("K132").Value = 4504.875
$ws.Range("L132").Value = 3157.0002
$ws.Range("M132").Value = -1974.875
$ws.Range("N132").Value = -8217.0002

$ws.Range("H136").Value = 4097.645
$ws.Range("I136").Value = 4054.5356
$ws.Range("J136").Value = 4500
$ws.Range("K136").Value = 12163.6068
$ws.Range("L136").Value = 13500
$ws.Range("M136").Value = -9613.606800000001
$ws.Range("N136").Value = -18600
